# Insert a new record row at row 236 in the "Pera" price sheet (weekly fruit/vegetable
# price update). This shifts all existing rows 236-334 down to 237-335, and the new
# row 236 carries a fresh weekly observation (new date + updated min/max/avg prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 236 - this shifts rows 236:334 down to 237:335,
# duplicating formatting from the row above (matches the observed style carry-over).
$ws.Rows.Item(236).Insert()

# Populate the newly inserted row 236 with the new weekly observation.
$ws.Cells.Item(236, 1).Value = 11
$ws.Cells.Item(236, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(236, 3).Value = "Bíobío"
$ws.Cells.Item(236, 4).Value = 44553
$ws.Cells.Item(236, 4).Style = "Normal"
$ws.Cells.Item(236, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(236, 5).Value = 8
$ws.Cells.Item(236, 6).Value = "Fruta"
$ws.Cells.Item(236, 7).Value = 100104
$ws.Cells.Item(236, 8).Value = "Frutos de pepita"
$ws.Cells.Item(236, 9).Value = 100104005
$ws.Cells.Item(236, 10).Value = "Pera"
$ws.Cells.Item(236, 11).Value = "Packham's Triumph"
$ws.Cells.Item(236, 12).Value = "Primera"
$ws.Cells.Item(236, 13).Value = 200
$ws.Cells.Item(236, 14).Value = 12000
$ws.Cells.Item(236, 15).Value = 13000
$ws.Cells.Item(236, 16).Value = 12500
$ws.Cells.Item(236, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(236, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(236, 19).Value = 781
$ws.Cells.Item(236, 20).Value = 16

# Append a duplicate of the (now shifted) final data row as the new last row 335,
# matching the trailing row in the edited workbook.
$ws.Cells.Item(335, 1).Value = 11
$ws.Cells.Item(335, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(335, 3).Value = "Bíobío"
$ws.Cells.Item(335, 4).Value = 44272
$ws.Cells.Item(335, 4).Style = "Normal"
$ws.Cells.Item(335, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(335, 5).Value = 8
$ws.Cells.Item(335, 6).Value = "Fruta"
$ws.Cells.Item(335, 7).Value = 100104
$ws.Cells.Item(335, 8).Value = "Frutos de pepita"
$ws.Cells.Item(335, 9).Value = 100104005
$ws.Cells.Item(335, 10).Value = "Pera"
$ws.Cells.Item(335, 11).Value = "Packham's Triumph"
$ws.Cells.Item(335, 12).Value = "Segunda"
$ws.Cells.Item(335, 13).Value = 100
$ws.Cells.Item(335, 14).Value = 8000
$ws.Cells.Item(335, 15).Value = 8000
$ws.Cells.Item(335, 16).Value = 8000
$ws.Cells.Item(335, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(335, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(335, 19).Value = 500
$ws.Cells.Item(335, 20).Value = 16
